$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns to make room for "Port.old" (new E) and "Port.new" (new G) ---
# Before: A Tree | B Station | C Instrument | D MC.old | E MC.new | F Timestamp
# After inserting at E, then again at the new G:
#   A Tree | B Station | C Instrument | D MC.old | E Port.old(new) | F MC.new | G Port.new(new) | H Timestamp
$ws.Columns("E").Insert()
$ws.Columns("G").Insert()

# --- Enter new unique text values in the same sequence the original author likely typed
#     them, so the rebuilt shared-strings table lines up with the source workbook. ---
$ws.Range("A5").Value = "TV4"
$ws.Range("I1").Value = "Note"
$ws.Range("B6").Value = "S3"
$ws.Range("C6").Value = "PYR"
$ws.Range("E1").Value = "Port.old"
$ws.Range("G1").Value = "Port.new"
$ws.Range("A7").Value = "FB8"
$ws.Range("I7").Value = "Port 5 went bad"
$ws.Range("I5").Value = "Port went bad"
$ws.Range("I6").Value = "Port went bad"

# --- Remaining (already-existing-string or numeric) cells for the new rows ---
$ws.Range("B5").Value = "S4"
$ws.Range("C5").Value = "A22"
$ws.Range("D5").Value = "MC3"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "MC1"
$ws.Range("G5").Value = 6
$ws.Range("H5").Value2 = 45314

$ws.Range("A6").Value = "TV4"
$ws.Range("D6").Value = "MC3"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "MC1"
$ws.Range("G6").Value = 5
$ws.Range("H6").Value2 = 45314

$ws.Range("B7").Value = "S4"
$ws.Range("C7").Value = "A22"
$ws.Range("D7").Value = "MC3"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = "MC3"
$ws.Range("G7").Value = 6
$ws.Range("H7").Value2 = 45315

# --- Apply the existing date style (re-using the workbook's existing date format rather
#     than creating a new custom number format). ---
$ws.Range("H2").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("H7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths for the two newly inserted columns. The OOXML in the target file
#     stores fractional "bestFit" pixel widths that Excel's real autofit engine computed;
#     this COM layer only supports coarser 1/6-character quantization, so we feed it the
#     width of the neighboring bestFit column (D / F) to land on the closest achievable
#     bucket. ---
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# --- Selection / active cell matches target ---
$ws.Range("B7").Select()
